# Apply the deck edit: append three new "Title and Content" slides
# (Key Shareholders / Valuation / Predictions) after the existing two slides.
#
# Slide 2 ("Technology") already uses the "Title and Content" layout that the
# new slides need, so each new slide is produced by duplicating the previous
# last slide (keeping the layout relationship + clrMapOvr intact) and then
# overwriting its title / body text in-place via the Characters() range (this
# replaces the text of the existing run instead of creating a brand-new run).

$p = $ppt.ActivePresentation

function Set-PlaceholderText {
    param($Shape, [string]$NewText)
    $tr = $Shape.TextFrame.TextRange
    $tr.Characters(1, $tr.Length).Text = $NewText
}

$slide2 = $p.Slides.Item(2)

$slide3 = $slide2.Duplicate().Item(1)
Set-PlaceholderText $slide3.Shapes.Item(1) "Key Shareholders"
Set-PlaceholderText $slide3.Shapes.Item(2) " Vanguard Group, Inc. (The)  &  Blackrock Inc. "

$slide4 = $slide3.Duplicate().Item(1)
Set-PlaceholderText $slide4.Shapes.Item(1) "Valuation"
Set-PlaceholderText $slide4.Shapes.Item(2) "Market Cap:  2.081T P/E Ratio:  37.32 EPS:  3.28 "

$slide5 = $slide4.Duplicate().Item(1)
Set-PlaceholderText $slide5.Shapes.Item(1) "Predictions"
Set-PlaceholderText $slide5.Shapes.Item(2) "Rating:  Overvalued Expected Return on Equity in 6 months based on TTM:  -14% Est. Return"
